$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 251.125
$ws.Range("I19").Value = 201.75
$ws.Range("K19").Value = 201.75
$ws.Range("M19").Value = -26.75
$ws.Range("H43").Value = 5566.3335
$ws.Range("I43").Value = 3812
$ws.Range("J43").Value = 6969.8
$ws.Range("K43").Value = 3812
$ws.Range("L43").Value = 6969.8
$ws.Range("M43").Value = -3743
$ws.Range("N43").Value = -7107.8
$ws.Range("H51").Value = 6832.8335
$ws.Range("I51").Value = 8749.25
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 8749.25
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -8265.25
$ws.Range("N51").Value = -3968
$ws.Range("H64").Value = 12117.294
$ws.Range("J64").Value = 19666.5
$ws.Range("L64").Value = 19666.5
$ws.Range("N64").Value = -20162.5
$ws.Range("H67").Value = 12117.294
$ws.Range("J67").Value = 19666.5
$ws.Range("L67").Value = 19666.5
$ws.Range("N67").Value = -21382.5
$ws.Range("H113").Value = 4910.8887
$ws.Range("I113").Value = 2799.6
$ws.Range("J113").Value = 7550
$ws.Range("K113").Value = 2799.6
$ws.Range("L113").Value = 7550
$ws.Range("M113").Value = 454.4000000000001
$ws.Range("N113").Value = -14058
$ws.Range("H116").Value = 6095.375
$ws.Range("I116").Value = 6256.2856
$ws.Range("J116").Value = 6029.1177
$ws.Range("K116").Value = 6256.2856
$ws.Range("L116").Value = 6029.1177
$ws.Range("M116").Value = -2814.2856
$ws.Range("N116").Value = -12913.1177
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 1618
$ws.Range("I12").Value = 1491
$ws.Range("J12").Value = 1999
$ws.Range("K12").Value = 1491
$ws.Range("L12").Value = 1999
$ws.Range("M12").Value = -1318
$ws.Range("N12").Value = -2345
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4788
$ws.Range("H74").Value = 2678.2
$ws.Range("I74").Value = 2678.2
$ws.Range("K74").Value = 2678.2
$ws.Range("M74").Value = -1804.2
$ws.Range("H77").Value = 2678.2
$ws.Range("I77").Value = 2678.2
$ws.Range("K77").Value = 13391
$ws.Range("M77").Value = -9023
$ws.Range("H82").Value = 60090.5
$ws.Range("J82").Value = 80181
$ws.Range("L82").Value = 80181
$ws.Range("N82").Value = -80903
$ws.Range("H85").Value = 60090.5
$ws.Range("J85").Value = 80181
$ws.Range("L85").Value = 80181
$ws.Range("N85").Value = -82677
$ws.Range("H132").Value = 750
$ws.Range("I132").Value = 750
$ws.Range("K132").Value = 2250
$ws.Range("M132").Value = 280
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3650.8
$ws.Range("I86").Value = 3249
$ws.Range("J86").Value = 4253.5
$ws.Range("K86").Value = 3249
$ws.Range("L86").Value = 4253.5
$ws.Range("M86").Value = -2126
$ws.Range("N86").Value = -6499.5
$ws.Range("H89").Value = 3650.8
$ws.Range("I89").Value = 3249
$ws.Range("J89").Value = 4253.5
$ws.Range("K89").Value = 16245
$ws.Range("L89").Value = 21267.5
$ws.Range("M89").Value = -10629
$ws.Range("N89").Value = -32499.5
$ws.Range("H105").Value = 3768.762
$ws.Range("I105").Value = 1581.6
$ws.Range("K105").Value = 1581.6
$ws.Range("M105").Value = 165.4000000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2985.4614
$ws.Range("I31").Value = 2859.4167
$ws.Range("J31").Value = 4498
$ws.Range("K31").Value = 2859.4167
$ws.Range("L31").Value = 4498
$ws.Range("M31").Value = -2564.4167
$ws.Range("N31").Value = -5088
$ws.Range("H34").Value = 2985.4614
$ws.Range("I34").Value = 2859.4167
$ws.Range("J34").Value = 4498
$ws.Range("K34").Value = 2859.4167
$ws.Range("L34").Value = 4498
$ws.Range("M34").Value = -2657.4167
$ws.Range("N34").Value = -4902
$ws.Range("J58").Value = 1000
$ws.Range("L58").Value = 1000
$ws.Range("N58").Value = -1406
$ws.Range("H105").Value = 1032.6471
$ws.Range("I105").Value = 1022.8125
$ws.Range("K105").Value = 1022.8125
$ws.Range("M105").Value = 724.1875
$ws.Range("H132").Value = 9048.700000000001
$ws.Range("I132").Value = 3097.5
$ws.Range("K132").Value = 9292.5
$ws.Range("M132").Value = -6762.5
$ws.Range("H134").Value = 2261.8386
$ws.Range("I134").Value = 2070.6333
$ws.Range("K134").Value = 6211.8999
$ws.Range("M134").Value = -3676.8999
$ws.Range("H135").Value = 50780
$ws.Range("J135").Value = 50780
$ws.Range("L135").Value = 50780
$ws.Range("N135").Value = -60920
$ws.Range("J136").Value = 1000
$ws.Range("L136").Value = 3000
$ws.Range("N136").Value = -8100
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 832.3333
$ws.Range("I75").Value = 799.25
$ws.Range("K75").Value = 2397.75
$ws.Range("M75").Value = -1399.75
$ws.Range("H78").Value = 832.3333
$ws.Range("I78").Value = 799.25
$ws.Range("K78").Value = 7193.25
$ws.Range("M78").Value = -2201.25
$ws.Range("H114").Value = 270.8889
$ws.Range("J114").Value = 49
$ws.Range("L114").Value = 147
$ws.Range("N114").Value = -6655
$ws.Range("H117").Value = 135
$ws.Range("J117").Value = 120
$ws.Range("L117").Value = 360
$ws.Range("N117").Value = -7244
$ws.Range("H121").Value = 488.75
$ws.Range("I121").Value = 492.33334
$ws.Range("J121").Value = 478
$ws.Range("K121").Value = 1477.00002
$ws.Range("L121").Value = 1434
$ws.Range("M121").Value = -167.0000199999999
$ws.Range("N121").Value = -4054
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3557.524
$ws.Range("I122").Value = 3461.353
$ws.Range("J122").Value = 3966.25
$ws.Range("K122").Value = 10384.059
$ws.Range("L122").Value = 11898.75
$ws.Range("M122").Value = -7934.059000000001
$ws.Range("N122").Value = -16798.75
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = ""
$ws.Range("H132").Value = 5673.6665
$ws.Range("I132").Value = 4010.5
$ws.Range("K132").Value = 12031.5
$ws.Range("M132").Value = -9501.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2415.8
$ws.Range("I61").Value = 1783.25
$ws.Range("K61").Value = 1783.25
$ws.Range("M61").Value = -1581.25
$ws.Range("H113").Value = 2415.8
$ws.Range("I113").Value = 1783.25
$ws.Range("K113").Value = 1783.25
$ws.Range("M113").Value = 386.75
$ws.Range("H132").Value = 4517.706
$ws.Range("I132").Value = 4237.5625
$ws.Range("K132").Value = 12712.6875
$ws.Range("M132").Value = -10182.6875
$ws.Range("H136").Value = 2174.75
$ws.Range("I136").Value = 2174.75
$ws.Range("K136").Value = 6524.25
$ws.Range("M136").Value = -3974.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2238.4
$ws.Range("I132").Value = 2104.3333
$ws.Range("K132").Value = 6312.999899999999
$ws.Range("M132").Value = -3782.999899999999
$ws.Range("H136").Value = 2261.0386
$ws.Range("I136").Value = 1686.6522
$ws.Range("K136").Value = 5059.9566
$ws.Range("M136").Value = -2509.9566
